# Update NATMI LR-pair output values (Efna1-Epha3) with new TPM-derived results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"23.630375"
$ws.Range("H2").Value = [double]"70.89112499999999"
$ws.Range("I2").Value = [double]"0.9002398112414131"
$ws.Range("J2").Value = [double]"0.9002398112414129"
$ws.Range("M2").Value = [double]"0.003058333333333333"
$ws.Range("N2").Value = [double]"0.009175000000000001"
$ws.Range("O2").Value = [double]"0.0001379486413073712"
$ws.Range("P2").Value = [double]"0.0001379486413073712"
$ws.Range("Q2").Value = [double]"0.07226956354166666"
$ws.Range("R2").Value = [double]"0.6504260718749999"
$ws.Range("S2").Value = [double]"0.0001241868588115573"
$ws.Range("T2").Value = [double]"0.0001241868588115572"
$ws.Range("G3").Value = [double]"23.630375"
$ws.Range("H3").Value = [double]"70.89112499999999"
$ws.Range("I3").Value = [double]"0.9002398112414131"
$ws.Range("J3").Value = [double]"0.9002398112414129"
$ws.Range("O3").Value = [double]"0.9939610820947024"
$ws.Range("P3").Value = [double]"0.9939610820947024"
$ws.Range("Q3").Value = [double]"520.7237483429166"
$ws.Range("R3").Value = [double]"4686.513735086249"
$ws.Range("S3").Value = [double]"0.8948033369262455"
$ws.Range("T3").Value = [double]"0.8948033369262454"
$ws.Range("G4").Value = [double]"23.630375"
$ws.Range("H4").Value = [double]"70.89112499999999"
$ws.Range("I4").Value = [double]"0.9002398112414131"
$ws.Range("J4").Value = [double]"0.9002398112414129"
$ws.Range("O4").Value = [double]"0.005900969263990248"
$ws.Range("P4").Value = [double]"0.005900969263990248"
$ws.Range("Q4").Value = [double]"3.091443809374999"
$ws.Range("R4").Value = [double]"27.822994284375"
$ws.Range("S4").Value = [double]"0.005312287456355961"
$ws.Range("T4").Value = [double]"0.005312287456355961"
$ws.Range("I5").Value = [double]"0.06214870537054815"
$ws.Range("J5").Value = [double]"0.06214870537054815"
$ws.Range("M5").Value = [double]"0.003058333333333333"
$ws.Range("N5").Value = [double]"0.009175000000000001"
$ws.Range("O5").Value = [double]"0.0001379486413073712"
$ws.Range("P5").Value = [double]"0.0001379486413073712"
$ws.Range("Q5").Value = [double]"0.004989181499999999"
$ws.Range("R5").Value = [double]"0.0449026335"
$ws.Range("S5").Value = [double]"8.573329464879243E-06"
$ws.Range("T5").Value = [double]"8.573329464879243E-06"
$ws.Range("I6").Value = [double]"0.06214870537054815"
$ws.Range("J6").Value = [double]"0.06214870537054815"
$ws.Range("O6").Value = [double]"0.9939610820947024"
$ws.Range("P6").Value = [double]"0.9939610820947024"
$ws.Range("S6").Value = [double]"0.06177339444089488"
$ws.Range("T6").Value = [double]"0.06177339444089488"
$ws.Range("I7").Value = [double]"0.06214870537054815"
$ws.Range("J7").Value = [double]"0.06214870537054815"
$ws.Range("O7").Value = [double]"0.005900969263990248"
$ws.Range("P7").Value = [double]"0.005900969263990248"
$ws.Range("S7").Value = [double]"0.0003667376001883903"
$ws.Range("T7").Value = [double]"0.0003667376001883903"
$ws.Range("I8").Value = [double]"0.03761148338803896"
$ws.Range("J8").Value = [double]"0.03761148338803896"
$ws.Range("M8").Value = [double]"0.003058333333333333"
$ws.Range("N8").Value = [double]"0.009175000000000001"
$ws.Range("O8").Value = [double]"0.0001379486413073712"
$ws.Range("P8").Value = [double]"0.0001379486413073712"
$ws.Range("Q8").Value = [double]"0.003019379341666667"
$ws.Range("R8").Value = [double]"0.027174414075"
$ws.Range("S8").Value = [double]"5.188453030934739E-06"
$ws.Range("T8").Value = [double]"5.188453030934739E-06"
$ws.Range("I9").Value = [double]"0.03761148338803896"
$ws.Range("J9").Value = [double]"0.03761148338803896"
$ws.Range("O9").Value = [double]"0.9939610820947024"
$ws.Range("P9").Value = [double]"0.9939610820947024"
$ws.Range("S9").Value = [double]"0.03738435072756213"
$ws.Range("T9").Value = [double]"0.03738435072756213"
$ws.Range("I10").Value = [double]"0.03761148338803896"
$ws.Range("J10").Value = [double]"0.03761148338803896"
$ws.Range("O10").Value = [double]"0.005900969263990248"
$ws.Range("P10").Value = [double]"0.005900969263990248"
$ws.Range("S10").Value = [double]"0.0002219442074458977"
$ws.Range("T10").Value = [double]"0.0002219442074458977"
